$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a value to be stored as text, guarding against Excel
# auto-converting numeric-looking strings (e.g. "21.20" -> 21.2, or
# "308.37" -> the number 308.37 instead of the text "308.37").
function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

$ws.Range('D2').Value = '44.159.71'
$ws.Range('E2').Value = '  +2.25%  '

$ws.Range('D3').Value = '2.429.76'
$ws.Range('E3').Value = '  +2.00%  '

$ws.Range('E4').Value = '  -0.04%  '

Set-TextValue 'D5' '308.37'
$ws.Range('E5').Value = '  +1.75%  '

Set-TextValue 'D6' '100.63'
$ws.Range('E6').Value = '  +3.71%  '

$ws.Range('E7').Value = '  +0.65%  '

$ws.Range('E9').Value = '  -0.33%  '

Set-TextValue 'D10' '35.44'
$ws.Range('E10').Value = '  +3.57%  '

Set-TextValue 'D11' '0.0801'
$ws.Range('E11').Value = '  +1.55%  '

$ws.Range('E12').Value = '  +2.76%  '

Set-TextValue 'D13' '18.73'
$ws.Range('E13').Value = '  +1.67%  '

$ws.Range('E14').Value = '  +2.16%  '

$ws.Range('D15').Value = '2.808.14'
$ws.Range('E15').Value = '  +1.96%  '

$ws.Range('D16').Value = '2.433.75'
$ws.Range('E16').Value = '  +2.74%  '

Set-TextValue 'D17' '0.833'
$ws.Range('E17').Value = '  +2.97%  '

$ws.Range('D18').Value = '44.124.96'
$ws.Range('E18').Value = '  +2.16%  '

$ws.Range('E19').Value = '  +1.21%  '

$ws.Range('E20').Value = '  +2.01%  '

$ws.Range('D21').Value = '0.0₃0907'
$ws.Range('E21').Value = '  +2.02%  '

Set-TextValue 'D22' '68.66'
$ws.Range('E22').Value = '  +0.07%  '

$ws.Range('B23').Value = 'BitcoinCash'
$ws.Range('C23').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue 'D23' '240.93'
$ws.Range('E23').Value = '  +2.45%  '

$ws.Range('B24').Value = 'ImmutableX'
$ws.Range('C24').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue 'D24' '2.30'
$ws.Range('E24').Value = '  +3.09%  '

$ws.Range('E25').Value = '  +1.77%  '

$ws.Range('E26').Value = '  -0.08%  '

$ws.Range('E27').Value = '  +1.74%  '

Set-TextValue 'D28' '2.34'
$ws.Range('E28').Value = '  -1.06%  '

Set-TextValue 'D29' '9.65'
$ws.Range('E29').Value = '  +5.70%  '

Set-TextValue 'D30' '33.25'
$ws.Range('E30').Value = '  +5.48%  '

$ws.Range('E31').Value = '  +15.41%  '

Set-TextValue 'D32' '18.71'
$ws.Range('E32').Value = '  +9.40%  '

Set-TextValue 'D33' '5.19'
$ws.Range('E33').Value = '  +1.72%  '

$ws.Range('E34').Value = '  +0.02%  '

Set-TextValue 'D35' '0.0763'
$ws.Range('E35').Value = '  +3.59%  '

$ws.Range('E36').Value = '  +3.53%  '

Set-TextValue 'D37' '4.55'
$ws.Range('E37').Value = '  +5.26%  '

Set-TextValue 'D38' '129.93'
$ws.Range('E38').Value = '  +22.98%  '

Set-TextValue 'D39' '2.93'
$ws.Range('E39').Value = '  +4.80%  '

$ws.Range('E40').Value = '  -0.63%  '

$ws.Range('E41').Value = '  +0.16%  '

Set-TextValue 'D42' '21.20'
$ws.Range('E42').Value = '  -5.03%  '

$ws.Range('E43').Value = '  +2.91%  '

$ws.Range('D44').Value = '1.963.84'
$ws.Range('E44').Value = '  +0.37%  '

$ws.Range('E45').Value = '  +1.78%  '

Set-TextValue 'D46' '2.89'
$ws.Range('E46').Value = '  +4.92%  '

Set-TextValue 'D47' '9.41'
$ws.Range('E47').Value = '  +1.39%  '

$ws.Range('E48').Value = '  +8.71%  '

$ws.Range('B49').Value = 'MultiversX'
$ws.Range('C49').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
Set-TextValue 'D49' '53.43'
$ws.Range('E49').Value = '  +1.19%  '

$ws.Range('B50').Value = 'BitcoinSV'
$ws.Range('C50').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
Set-TextValue 'D50' '73.60'
$ws.Range('E50').Value = '  +2.37%  '

$ws.Range('B51').Value = 'TrustWalletToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue 'D51' '1.16'
$ws.Range('E51').Value = '  +1.42%  '
